$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.516.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.15%  "

# Row 3
$ws.Range("D3").Value = "'1.672.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.68%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "'219.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.23%  "

# Row 6
$ws.Range("D6").Value = "'0.5283"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.12%  "

# Row 7
$ws.Range("D7").Value = "'1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").Value = "'0.2682"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.97%  "

# Row 9
$ws.Range("D9").Value = "'0.06385"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.32%  "

# Row 10
$ws.Range("D10").Value = "'21.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.84%  "

# Row 11
$ws.Range("D11").Value = "'0.07804"
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'1.676.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.99%  "

# Row 13
$ws.Range("D13").Value = "'4.487"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.39%  "

# Row 14
$ws.Range("D14").Value = "'0.5577"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.23%  "

# Row 15
$ws.Range("D15").Value = "'0.0₅8320"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "

# Row 16
$ws.Range("D16").Value = "'65.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.70%  "

# Row 17
$ws.Range("D17").Value = "'26.525.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.12%  "

# Row 18
$ws.Range("E18").Value = "  -0.01%  "

# Row 19
$ws.Range("D19").Value = "'4.762"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.36%  "

# Row 20
$ws.Range("D20").Value = "'193.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.27%  "

# Row 21
$ws.Range("D21").Value = "'10.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.96%  "

# Row 22
$ws.Range("D22").Value = "'6.319"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.99%  "

# Row 23
$ws.Range("D23").Value = "'1.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.12%  "

# Row 24
$ws.Range("D24").Value = "'0.1267"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.45%  "

# Row 25
$ws.Range("D25").Value = "'139.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.34%  "

# Row 26
$ws.Range("D26").Value = "'7.413"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.28%  "

# Row 27
$ws.Range("D27").Value = "'16.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.64%  "

# Row 28
$ws.Range("D28").Value = "'1.424"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.15%  "

# Row 29
$ws.Range("D29").Value = "'0.06199"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.10%  "

# Row 30
$ws.Range("D30").Value = "'1.292"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.29%  "

# Row 31
$ws.Range("D31").Value = "'3.620"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.54%  "

# Row 32
$ws.Range("D32").Value = "'3.424"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.04%  "

# Row 33
$ws.Range("D33").Value = "'1.683"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.44%  "

# Row 34
$ws.Range("D34").Value = "'1.007"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.63%  "

# Row 35
$ws.Range("D35").Value = "'0.6084"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.02%  "

# Row 36
$ws.Range("D36").Value = "'2.414"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.90%  "

# Row 37
$ws.Range("D37").Value = "'2.779"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.08%  "

# Row 38
$ws.Range("D38").Value = "'0.01616"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.75%  "

# Row 39
$ws.Range("D39").Value = "'6.072"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.33%  "

# Row 40
$ws.Range("D40").Value = "'1.094.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.54%  "

# Row 41
$ws.Range("D41").Value = "'0.8565"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.41%  "

# Row 42
$ws.Range("D42").Value = "'1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("E43").Value = "  +2.05%  "

# Row 44
$ws.Range("D44").Value = "'1.818.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.44%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'58.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.67%  "

# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.0₈108"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.20%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'8.143"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.10%  "

# Row 48
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.004"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.39%  "

# Row 49
$ws.Range("D49").Value = "'1.514"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.25%  "

# Row 50
$ws.Range("D50").Value = "'0.05208"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.20%  "

# Row 51
$ws.Range("D51").Value = "'6.004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.08%  "

